# Update "想去人数" (column F) counts on the "展览" and "全部类型" sheets.
# Each row is identified by its event name (column C) *and* its current F
# value, since a couple of event names repeat across different rows/dates
# with different counts (e.g. "杭州·TCD国潮动漫游戏嘉年华" appears twice on
# the "展览" sheet). Matching on (name, old value) avoids touching the
# wrong occurrence.

$wb = $excel.ActiveWorkbook

function Apply-Updates($sheetName, $changes) {
    $ws = $wb.Worksheets.Item($sheetName)
    $usedRange = $ws.UsedRange
    $rowCount = $usedRange.Rows.Count

    for ($r = 2; $r -le $rowCount; $r++) {
        $name = $ws.Cells.Item($r, 3).Value()
        $current = $ws.Cells.Item($r, 6).Value()
        foreach ($change in $changes) {
            if ($name -eq $change.Name -and $current -eq $change.Old) {
                $ws.Cells.Item($r, 6).Value = $change.New
                break
            }
        }
    }
}

$changes = @(
    @{ Name = "杭州·CICAF·中国国风品牌盛典、中国COSPLAY超级盛典"; Old = 640; New = 642 }
    @{ Name = "杭州·第二十届中国国际动漫节主会场门票"; Old = 6068; New = 6085 }
    @{ Name = "杭州·蔚蓝档案only"; Old = 463; New = 464 }
    @{ Name = "杭州·热血番&运动番ONLY"; Old = 408; New = 410 }
    @{ Name = "杭州·第38届漫展x原崩铁only"; Old = 1385; New = 1388 }
    @{ Name = "杭州·第三届日夜国乙only"; Old = 3137; New = 3141 }
    @{ Name = "杭州·代号鸢沧笙踏歌only(免票)"; Old = 404; New = 439 }
    @{ Name = "杭州·次元盛典1.0"; Old = 1983; New = 1988 }
    @{ Name = "杭州·第二届白日梦次元动漫嘉年华"; Old = 90; New = 91 }
    @{ Name = "杭州·乌托邦次元聚会3.0·二次元全女性夜场"; Old = 1007; New = 1008 }
    @{ Name = "杭州·黑执事only"; Old = 93; New = 102 }
    @{ Name = "杭州·AD04动漫展"; Old = 3748; New = 3758 }
    @{ Name = "杭州·代号鸢only-广陵大学"; Old = 1187; New = 1189 }
    @{ Name = "杭州·草莓动漫节"; Old = 2978; New = 2982 }
    @{ Name = "【会员购严选】杭州·首届次元格子动漫展-进入格子空间，探索次元世界！"; Old = 2532; New = 2533 }
    @{ Name = "杭州·TCD国潮动漫游戏嘉年华"; Old = 4320; New = 4327 }
    @{ Name = "杭州·生如夏花国乙only·日夜场"; Old = 483; New = 486 }
    @{ Name = "杭州·第五届华盟次元嘉年华&周年庆狂欢"; Old = 1365; New = 1367 }
    @{ Name = "杭州·第四届ArknightsOnly·狼与黑荆棘（明日方舟Only）"; Old = 136; New = 139 }
    @{ Name = "杭州·ACG CLUB动漫游戏嘉年华"; Old = 36; New = 37 }
    @{ Name = "杭州·夏之誓国乙only-日夜场"; Old = 1044; New = 1046 }
    @{ Name = "杭州·AP动漫游戏嘉年华"; Old = 1130; New = 1132 }
    @{ Name = "杭州·梦漫星河动漫展"; Old = 732; New = 734 }
    @{ Name = "杭州·梦漫星河动漫嘉年华·赵路专场"; Old = 648; New = 649 }
    @{ Name = "杭州·原神X星铁X绝区零only"; Old = 441; New = 442 }
    @{ Name = "杭州·造梦探险家城堡二次元同好会"; Old = 25; New = 26 }
    @{ Name = "杭州·HD·01"; Old = 129; New = 131 }
    @{ Name = "浙江·蔚蓝档案ONLY02-夏末狂欢！"; Old = 7; New = 8 }
    @{ Name = "杭州·D3动漫游戏嘉年华"; Old = 325; New = 326 }
    @{ Name = "杭州·理想乡动漫展-同人创作者大会"; Old = 3624; New = 3625 }
)

Apply-Updates "展览" $changes
Apply-Updates "全部类型" $changes
